$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Extra mobile data requests")

# Update the phone number in B3 (normalised by phonelib during save)
$ws.Range("B3").Value = "07900222222"

# Update the active selection to match the edited cell
$ws.Range("B3").Select()
